# Fixed skipped counting issue
# Update histogram bin counts in column B (Sheet1) to corrected values.
# The chart on the sheet references Sheet1!$B$1:$B$41, so updating the
# cell values will refresh the chart's cached series values as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    4  = 1
    6  = 1
    7  = 45
    8  = 92
    9  = 96
    10 = 60
    11 = 96
    12 = 59
    13 = 59
    14 = 25
    15 = 27
    16 = 61
    17 = 105
    18 = 191
    19 = 179
    20 = 234
    21 = 370
    22 = 555
    23 = 581
    24 = 519
    25 = 349
    26 = 285
    27 = 223
    28 = 181
    29 = 159
    30 = 109
    31 = 43
    32 = 9
    33 = 17
    34 = 12
    35 = 78
    36 = 219
    37 = 405
    38 = 1005
    39 = 2075
    40 = 3963
    41 = 31328
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}

$excel.CalculateFullRebuild()
